$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two new rows at the top of the data table (rows 2:3), pushing
# the existing data (and the totals row) down by two rows.
$ws.Rows("2:3").Insert()

# New daily data for 27 May (row 2) and 26 May (row 3). Grab the date
# number format from an existing date cell so A2:A3 match the rest of
# the column.
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

$ws.Range("A2").Value = 45439
$ws.Range("B2").Value = 253
$ws.Range("C2").Value = 32
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 217

$ws.Range("A3").Value = 45438
$ws.Range("B3").Value = 186
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 162

# Restore the totals-row formulas to cover the full, now-larger range.
$ws.Range("B29").Formula = "=SUM(B2:B28)"
$ws.Range("C29").Formula = "=SUM(C2:C28)"
$ws.Range("D29").Formula = "=SUM(D2:D28)"
$ws.Range("E29").Formula = "=SUM(E2:E28)"

# Match the author's final view state: scrolled down a couple of rows,
# with G28 selected (was G24 before the two rows were inserted).
$ws.Activate() | Out-Null
$ws.Range("G28").Select() | Out-Null
